# Generate Report for Archive
#
# The 41966e5f-3990-417d-a9c8-797c9abc00bd file moved up in the status
# report (now "In Translation" instead of "Ready for handoff"), so its
# row is relocated to directly follow the first (ef4783de) row, pushing
# the 6b9ff258 / 965104be / 253d137b rows down by one position each, on
# all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "41966e5f-3990-417d-a9c8-797c9abc00bd.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "2016-25-11 16:25:49"

$ws.Range("A4").Value = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "2016-25-11 16:25:17"

$ws.Range("A5").Value = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "2016-25-11 16:25:17"

$ws.Range("A6").Value = "253d137b-9592-410f-9fca-d89327456d1f.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "2016-24-11 16:24:14"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "41966e5f-3990-417d-a9c8-797c9abc00bd.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "41966e5f-3990-417d-a9c8-797c9abc00bd.570ead55e97107530a4552f7d5383995dba53976.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-11 16:25:46"

$ws.Range("A4").Value = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-11 16:25:14"

$ws.Range("A5").Value = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-11 16:25:14"

$ws.Range("A6").Value = "253d137b-9592-410f-9fca-d89327456d1f.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-11 16:24:11"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "41966e5f-3990-417d-a9c8-797c9abc00bd.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("D3").Value = "41966e5f-3990-417d-a9c8-797c9abc00bd.570ead55e97107530a4552f7d5383995dba53976.de-de.xlf"
$ws.Range("E3").Value = "2016-03-11 16:25:49"

$ws.Range("A4").Value = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("D4").Value = "6b9ff258-5cf7-40c6-93dc-10784c1d9a46.0abd8208a0eb44c32d9f52395849381cb7361d5f.de-de.xlf"
$ws.Range("E4").Value = "2016-03-11 16:25:17"

$ws.Range("A5").Value = "965104be-0ef5-4edb-82b1-facb0f37c968.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "965104be-0ef5-4edb-82b1-facb0f37c968.a05aeaf165ae5501d28d8aece37eedefb60075ad.de-de.xlf"
$ws.Range("E5").Value = "2016-03-11 16:25:17"

$ws.Range("A6").Value = "253d137b-9592-410f-9fca-d89327456d1f.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "253d137b-9592-410f-9fca-d89327456d1f.fc74a366d75ceb48353102f324d1f380db212c73.de-de.xlf"
$ws.Range("E6").Value = "2016-03-11 16:24:14"
